$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data refresh: update regcntr_id (column A) values from 10002/10005 -> 10003
# for the specified rows, per the 2nd May data refresh.
$rowsToUpdate = @(3, 23, 43, 63, 83, 105, 114, 123, 132, 141)
foreach ($r in $rowsToUpdate) {
    $ws.Cells.Item($r, 1).Value = 10003
}

# Update the sheet view: scroll position reset (no frozen topLeftCell) and
# selection moved to the first empty row below the data table (A162),
# selecting the full rows beneath the data (A162:XFD1048576).
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("A162:XFD1048576").Select()
